# Append: 2025-09-20 12:42 JST
# Update the "取得日時" (retrieved-at) timestamp in column A for all
# existing data rows on the scraping results sheet ("ランサーズ") to
# reflect the latest run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-20 12:42:15"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value -ne $null -and $ws.Cells.Item($r, 1).Value -ne "") {
        $ws.Cells.Item($r, 1).Value = $newTimestamp
    }
}
